$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New summary rows at the bottom (14-17): label in column A, aggregate formula in column B ---
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"

$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Style B14 (bold, size 12, vertically centered), then propagate the same format to B15:B17
$b14 = $ws.Range("B14")
$b14.Font.Bold = $true
$b14.Font.Size = 12
$b14.VerticalAlignment = -4108

$b14.Copy()
$ws.Range("B15:B17").PasteSpecial(-4122)

# --- Column-average cell at J12 (bold) ---
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"
$ws.Range("J12").Font.Bold = $true

# --- Page setup (printed as portrait on the standard small paper size) ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Leave the selection where the author left it
$ws.Range("J12").Select()
